# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Fri Aug 25 21:24:32 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.122.24'
$ws.Range('E2').Value = '  -0.29%  '

# Row 3
$ws.Range('D3').Value = '1.653.54'
$ws.Range('E3').Value = '  -0.42%  '

# Row 5
$ws.Range('E5').Value = '  +0.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5288'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.51%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.27%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2606'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.13%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06316'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.14%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.39'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.98%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07753'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.61%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.480'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.25%  '

# Row 13
$ws.Range('D13').Value = '1.660.39'
$ws.Range('E13').Value = '  +0.18%  '

# Row 14
$ws.Range('E14').Value = '  +0.02%  '

# Row 15
$ws.Range('D15').Value = '0.0₅8128'
$ws.Range('E15').Value = '  -0.86%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.86%  '

# Row 17
$ws.Range('D17').Value = '26.138.10'
$ws.Range('E17').Value = '  -0.35%  '

# Row 18
$ws.Range('E18').Value = '  -0.32%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.538'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.60%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.86%  '

# Row 21
$ws.Range('E21').Value = '  -1.02%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.981'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.51%  '

# Row 23
$ws.Range('E23').Value = '  -0.41%  '

# Row 25
$ws.Range('E25').Value = '  +0.36%  '

# Row 26
$ws.Range('E26').Value = '  +0.71%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.11%  '

# Row 28
$ws.Range('E28').Value = '  +1.59%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05933'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.09%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.277'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.31%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.502'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.40%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.228'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.59%  '

# Row 33
$ws.Range('E33').Value = '  -5.63%  '

# Row 34
$ws.Range('E34').Value = '  -0.06%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9441'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.43%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.759'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.77%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5630'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.20%  '

# Row 38
$ws.Range('E38').Value = '  +1.11%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.847'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.59%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8442'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.08%  '

# Row 41
$ws.Range('E41').Value = '  -0.19%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.31%  '

# Row 43
$ws.Range('D43').Value = '1.007.96'
$ws.Range('E43').Value = '  -2.35%  '

# Row 44
$ws.Range('D44').Value = '1.798.82'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '56.81'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.29%  '

# Row 46
$ws.Range('E46').Value = '  -1.05%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.56%  '

# Row 48
$ws.Range('E48').Value = '  +1.37%  '

# Row 49
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05151'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.65%  '

# Row 50
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.471'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.60%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.740'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.34%  '
